$wb = $excel.ActiveWorkbook

# --- Adform sheet: move selection, it is no longer the active tab ---
$wsAdform = $wb.Worksheets.Item("Adform")
$wsAdform.Range("C24").Select()

# --- TTD sheet: add a new "Price Type" column and a new rate row ---
$wsTTD = $wb.Worksheets.Item("TTD")

# Insert a new column before the existing "Batch ID" column (I) for "Price Type"
$wsTTD.Columns("I:I").Insert()

# Clear the leftover formatting on the old placeholder row 4 cells so the
# new data row carries no stray styles
$wsTTD.Range("H4:J4").ClearFormats()

# Match the column width used by the sibling "Price" column
$wsTTD.Columns("I:I").ColumnWidth = 20.666666666666668

# Header + description for the new column
$wsTTD.Range("I1").Value = "Price Type"
$wsTTD.Range("I2").Value = "Add/Edit: Required" + [char]10 + "Edit Rates: Required" + [char]10 + "Retrieve Batch: Optional" + [char]10 + "Retrieve Rates: Optional" + [char]10 + "Values: CPM or PercentOfMediaCost"

# The existing rate row now needs an explicit Price Type value
$wsTTD.Range("I3").Value = "CPM"

# New rate row showing a percent-of-media-cost rate
$wsTTD.Range("A4").Value = 20190401004
$wsTTD.Range("B4").Value = "ttdratetest_partnerID_rate"
$wsTTD.Range("E4").Value = $true
$wsTTD.Range("H4").Value = 1
$wsTTD.Range("I4").Value = "PercentOfMediaCost"
$wsTTD.Range("F4").Value = "eyeota"
$wsTTD.Range("G4").Value = "abc123"
$wsTTD.Range("C4").Value = "Test Segment 20190401004"
$wsTTD.Range("D4").Value = "Test Segment 20190401004"
$wsTTD.Range("J4").Value = 23456

# Description row grows to fit the extra column's wrapped text
$wsTTD.Rows("2:2").RowHeight = 102

# TTD becomes the active sheet/tab
$wsTTD.Activate()
$wsTTD.Range("B5").Select()
